$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("A2").Value = 58191
$ws.Range("B2").Value = "Ana Clara Oliveira"
$ws.Range("C2").Value = "Atendimento ao Cliente"
$ws.Range("D2").Value = "Outros"
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 45086
$ws.Range("G2").Value = 9020.02

# Row 3
$ws.Range("A3").Value = 92980
$ws.Range("B3").Value = "Dr. Vinicius Jesus"
$ws.Range("C3").Value = "Vendas"
$ws.Range("D3").Value = "Consulta médica"
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = 45098
$ws.Range("G3").Value = 4578.26

# Row 4
$ws.Range("A4").Value = 28691
$ws.Range("B4").Value = "Ana Lívia Alves"
$ws.Range("C4").Value = "Atendimento ao Cliente"
$ws.Range("D4").Value = "Viagem de negócios"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 45088
$ws.Range("G4").Value = 5435.91

# Row 5
$ws.Range("A5").Value = 11570
$ws.Range("B5").Value = "Vitor Hugo Souza"
$ws.Range("C5").Value = "Financeiro"
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 45103
$ws.Range("G5").Value = 7685.4

# Row 6
$ws.Range("A6").Value = 9880
$ws.Range("B6").Value = "Daniela Oliveira"
$ws.Range("C6").Value = "Vendas"
$ws.Range("D6").Value = "Problemas pessoais"
$ws.Range("E6").Value = 6
$ws.Range("F6").Value = 45086
$ws.Range("G6").Value = 4199.66

# Row 7
$ws.Range("A7").Value = 15574
$ws.Range("B7").Value = "Maria Luiza Vieira"
$ws.Range("C7").Value = "Marketing"
$ws.Range("D7").Value = "Consulta médica"
$ws.Range("E7").Value = 8
$ws.Range("F7").Value = 45086
$ws.Range("G7").Value = 4049.41

# Row 8
$ws.Range("A8").Value = 23628
$ws.Range("B8").Value = "Stephany Moraes"
$ws.Range("C8").Value = "Recursos Humanos"
$ws.Range("D8").Value = "Viagem de negócios"
$ws.Range("E8").Value = 8
$ws.Range("F8").Value = 45097
$ws.Range("G8").Value = 10784.02

# Row 9
$ws.Range("A9").Value = 4365
$ws.Range("B9").Value = "Maria Fernanda Costela"
$ws.Range("C9").Value = "Recursos Humanos"
$ws.Range("D9").Value = "Outros"
$ws.Range("E9").Value = 8
$ws.Range("F9").Value = 45078
$ws.Range("G9").Value = 4575.23

# Row 10
$ws.Range("A10").Value = 95373
$ws.Range("B10").Value = "Isabelly Silveira"
$ws.Range("C10").Value = "Marketing"
$ws.Range("D10").Value = "Viagem de negócios"
$ws.Range("E10").Value = 7
$ws.Range("F10").Value = 45092
$ws.Range("G10").Value = 8054.53

# Row 11
$ws.Range("A11").Value = 9757
$ws.Range("B11").Value = "Dr. Francisco Lopes"
$ws.Range("C11").Value = "P&D"
$ws.Range("D11").Value = "Consulta médica"
$ws.Range("E11").Value = 8
$ws.Range("F11").Value = 45093
$ws.Range("G11").Value = 8084.98
